# Helper: EMU -> points, with a small epsilon nudge so that float round-trip
# through the host's point-based COM properties lands back on the exact EMU
# value instead of truncating one unit short. Zero stays exactly zero so the
# "minimum 1 EMU" floor some shape constructors apply never kicks in.
function EMU($emu) {
    if ($emu -eq 0) { return 0 }
    return ($emu / 12700.0) + 0.00003
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: wrap the two pictures + two caption textboxes in a new group,
# and widen/reflow the right-hand ("Philip Wolfe") caption textbox.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$range = $s1.Shapes.Range(@(1, 2, 3, 4))
$grp11 = $range.Group()
$grp11.Name = "组合 11"

$wolfeBox = $grp11.GroupItems.Item(4)
$wolfeBox.TextFrame.WordWrap = -1
$wolfeBox.Left = EMU 8276623
$wolfeBox.Top = EMU 6211669
$wolfeBox.Width = EMU 2496480
$wolfeBox.Height = EMU 923330

# ---------------------------------------------------------------------
# Slide 2 (new): right-triangle / coordinate-axes diagram.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 7)

# --- group "组合 10": three axis lines -------------------------------
$a = $s2.Shapes.AddLine(0, 0, 10, 10)
$a.Left = EMU 704193
$a.Top = EMU 1681655
$a.Width = EMU 840828
$a.Height = EMU 2511973

$b = $s2.Shapes.AddLine(0, 0, 10, 10)
$b.Left = EMU 1545021
$b.Top = EMU 4193628
$b.Width = EMU 3626069
$b.Height = EMU 0

$c = $s2.Shapes.AddLine(0, 0, 10, 10)
$c.Left = EMU 5171090
$c.Top = EMU 1681655
$c.Width = EMU 0
$c.Height = EMU 2511973
$c.VerticalFlip = $true

foreach ($ln in @($a, $b, $c)) {
    $ln.Line.Weight = EMU 15875
    $ln.Line.ForeColor.RGB = 0
}
$a.Name = "直接连接符 5"
$b.Name = "直接连接符 7"
$c.Name = "直接连接符 9"

$grp10 = $s2.Shapes.Range(@(1, 2, 3)).Group()
$grp10.Name = "组合 10"
$grp10.Left = EMU 2969873
$grp10.Top = EMU 1874695

# --- red dashed-drop arrow connectors --------------------------------
$arrow1 = $s2.Shapes.AddConnector(1, 0, 0, 10, 10)
$arrow1.Left = EMU 3810701
$arrow1.Top = EMU 1127760
$arrow1.Width = EMU 0
$arrow1.Height = EMU 3258908
$arrow1.VerticalFlip = $true
$arrow1.Line.Weight = EMU 12700
$arrow1.Line.ForeColor.RGB = 255
$arrow1.Line.EndArrowheadStyle = 2
$arrow1.Name = "直接箭头连接符 22"

$arrow2 = $s2.Shapes.AddConnector(1, 0, 0, 10, 10)
$arrow2.Left = EMU 2702560
$arrow2.Top = EMU 1127760
$arrow2.Width = EMU 1108141
$arrow2.Height = EMU 3258908
$arrow2.HorizontalFlip = $true
$arrow2.VerticalFlip = $true
$arrow2.Line.Weight = EMU 12700
$arrow2.Line.ForeColor.RGB = 255
$arrow2.Line.EndArrowheadStyle = 2
$arrow2.Name = "直接箭头连接符 24"

# --- group "组合 39": small tick mark ---------------------------------
$t1a = $s2.Shapes.AddLine(0, 0, 10, 10)
$t1a.Left = EMU 3810701
$t1a.Top = EMU 4145280
$t1a.Width = EMU 263459
$t1a.Height = EMU 0

$t1b = $s2.Shapes.AddLine(0, 0, 10, 10)
$t1b.Left = EMU 4074160
$t1b.Top = EMU 4155440
$t1b.Width = EMU 0
$t1b.Height = EMU 210908

foreach ($ln in @($t1a, $t1b)) {
    $ln.Line.Weight = EMU 9525
    $ln.Line.ForeColor.RGB = 0
}
$t1a.Name = "直接连接符 36"
$t1b.Name = "直接连接符 38"

$tickIdx = $s2.Shapes.Count
$grp39 = $s2.Shapes.Range(@($tickIdx - 1, $tickIdx)).Group()
$grp39.Name = "组合 39"
$grp39.Left = EMU 3810702
$grp39.Top = EMU 4216400
$grp39.Width = EMU 111058
$grp39.Height = EMU 149947

# --- group "组合 40": second small tick mark, rotated -----------------
$t2a = $s2.Shapes.AddLine(0, 0, 10, 10)
$t2a.Left = EMU 3810701
$t2a.Top = EMU 4145280
$t2a.Width = EMU 263459
$t2a.Height = EMU 0

$t2b = $s2.Shapes.AddLine(0, 0, 10, 10)
$t2b.Left = EMU 4074160
$t2b.Top = EMU 4155440
$t2b.Width = EMU 0
$t2b.Height = EMU 210908

foreach ($ln in @($t2a, $t2b)) {
    $ln.Line.Weight = EMU 9525
    $ln.Line.ForeColor.RGB = 0
}
$t2a.Name = "直接连接符 41"
$t2b.Name = "直接连接符 42"

$tick2Idx = $s2.Shapes.Count
$grp40 = $s2.Shapes.Range(@($tick2Idx - 1, $tick2Idx)).Group()
$grp40.Name = "组合 40"
$grp40.Left = EMU 3810702
$grp40.Top = EMU 4216400
$grp40.Width = EMU 111058
$grp40.Height = EMU 149947
$grp40.Rotation = 270
$grp40.Left = EMU 7290501
$grp40.Top = EMU 4237595

# --- right-triangle shape ---------------------------------------------
$tri = $s2.Shapes.AddShape(8, 0, 0, 10, 10)
$tri.Left = EMU 2875279
$tri.Top = EMU 1595119
$tri.Width = EMU 935419
$tri.Height = EMU 2771224
$tri.Rotation = 180
$tri.Name = "直角三角形 43"
$tri.Fill.ForeColor.ObjectThemeColor = 5
$tri.Line.ForeColor.ObjectThemeColor = 5
try { $tri.TextFrame.VerticalAnchor = 3 } catch {}
try { $tri.TextFrame.TextRange.ParagraphFormat.Alignment = 2 } catch {}
